# "derived volt creation ui done"
# Adds a new config row (row 7) for the "derived voltage creation service"
# URL, mirroring the existing raw-voltage / raw-frequency / ... rows above
# it: a label in column A and a hyperlinked URL (same google.com URL used
# by the sibling rows) in column B with the "Hyperlink" style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label / value pair on row 7.
$ws.Range("A7").Value = "derivedVoltageCreationServiceUrl"
$ws.Range("B7").Value = "http://google.com"

# Hyperlink B7 to the same external URL the other rows use.
$ws.Hyperlinks.Add($ws.Range("B7"), "http://google.com/") | Out-Null

# Match the visual style ("Hyperlink") already used by B3:B6.
$ws.Range("B7").Style = "Hyperlink"

# Leave the selection where the author left it after adding the row.
$ws.Range("B9").Select() | Out-Null
